$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.052.04"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.163.94"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "604.98"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "153.60"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.159.94"
$ws.Range("E8").Value = "  -1.59%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.545"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  -1.28%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -8.58%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  -1.34%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "38.24"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "3.683.61"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "66.086.47"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "3.164.77"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("E19").Value = "  +0.95%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "508.59"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "15.36"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("E24").Value = "  -4.05%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "84.46"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  +0.05%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.99"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.09"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.37"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +4.76%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.01"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.55%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +4.95%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "27.86"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -1.58%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "498.54"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.94%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "55.03"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0877"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.81%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0418"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  +6.14%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -4.81%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "2.814.75"
$ws.Range("E46").Value = "  -4.55%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "27.91"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("E50").Value = "  +0.44%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "35.20"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +6.18%  "
